$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare column B (text) cells so numeric-looking codes stay as text, not numbers
$bRange = $ws.Range("B114:B158")
$bRange.NumberFormat = "@"

$ws.Cells.Item(114,1).Value2 = 112
$ws.Cells.Item(114,2).Value = "215"
$ws.Cells.Item(114,3).Value2 = 321
$ws.Cells.Item(115,1).Value2 = 113
$ws.Cells.Item(115,2).Value = "260"
$ws.Cells.Item(115,3).Value2 = 321
$ws.Cells.Item(116,1).Value2 = 114
$ws.Cells.Item(116,2).Value = "268"
$ws.Cells.Item(116,3).Value2 = 321
$ws.Cells.Item(117,1).Value2 = 115
$ws.Cells.Item(117,2).Value = "288"
$ws.Cells.Item(117,3).Value2 = 321
$ws.Cells.Item(118,1).Value2 = 116
$ws.Cells.Item(118,2).Value = "355/1"
$ws.Cells.Item(118,3).Value2 = 321
$ws.Cells.Item(119,1).Value2 = 117
$ws.Cells.Item(119,2).Value = "356/1"
$ws.Cells.Item(119,3).Value2 = 321
$ws.Cells.Item(120,1).Value2 = 118
$ws.Cells.Item(120,2).Value = "356/2"
$ws.Cells.Item(120,3).Value2 = 321
$ws.Cells.Item(121,1).Value2 = 119
$ws.Cells.Item(121,2).Value = "356/5"
$ws.Cells.Item(121,3).Value2 = 321
$ws.Cells.Item(122,1).Value2 = 120
$ws.Cells.Item(122,2).Value = "356/25"
$ws.Cells.Item(122,3).Value2 = 321
$ws.Cells.Item(123,1).Value2 = 121
$ws.Cells.Item(123,2).Value = "356/27"
$ws.Cells.Item(123,3).Value2 = 321
$ws.Cells.Item(124,1).Value2 = 122
$ws.Cells.Item(124,2).Value = "356/28"
$ws.Cells.Item(124,3).Value2 = 321
$ws.Cells.Item(125,1).Value2 = 123
$ws.Cells.Item(125,2).Value = "356/29"
$ws.Cells.Item(125,3).Value2 = 321
$ws.Cells.Item(126,1).Value2 = 124
$ws.Cells.Item(126,2).Value = "357"
$ws.Cells.Item(126,3).Value2 = 321
$ws.Cells.Item(127,1).Value2 = 125
$ws.Cells.Item(127,2).Value = "358"
$ws.Cells.Item(127,3).Value2 = 321
$ws.Cells.Item(128,1).Value2 = 126
$ws.Cells.Item(128,2).Value = "361/1"
$ws.Cells.Item(128,3).Value2 = 321
$ws.Cells.Item(129,1).Value2 = 127
$ws.Cells.Item(129,2).Value = "361/2"
$ws.Cells.Item(129,3).Value2 = 321
$ws.Cells.Item(130,1).Value2 = 128
$ws.Cells.Item(130,2).Value = "362/1"
$ws.Cells.Item(130,3).Value2 = 321
$ws.Cells.Item(131,1).Value2 = 129
$ws.Cells.Item(131,2).Value = "362/31"
$ws.Cells.Item(131,3).Value2 = 321
$ws.Cells.Item(132,1).Value2 = 130
$ws.Cells.Item(132,2).Value = "362/33"
$ws.Cells.Item(132,3).Value2 = 321
$ws.Cells.Item(133,1).Value2 = 131
$ws.Cells.Item(133,2).Value = "362/34"
$ws.Cells.Item(133,3).Value2 = 321
$ws.Cells.Item(134,1).Value2 = 132
$ws.Cells.Item(134,2).Value = "362/35"
$ws.Cells.Item(134,3).Value2 = 321
$ws.Cells.Item(135,1).Value2 = 133
$ws.Cells.Item(135,2).Value = "362/36"
$ws.Cells.Item(135,3).Value2 = 321
$ws.Cells.Item(136,1).Value2 = 134
$ws.Cells.Item(136,2).Value = "363"
$ws.Cells.Item(136,3).Value2 = 321
$ws.Cells.Item(137,1).Value2 = 135
$ws.Cells.Item(137,2).Value = "364"
$ws.Cells.Item(137,3).Value2 = 321
$ws.Cells.Item(138,1).Value2 = 136
$ws.Cells.Item(138,2).Value = "388"
$ws.Cells.Item(138,3).Value2 = 321
$ws.Cells.Item(139,1).Value2 = 137
$ws.Cells.Item(139,2).Value = "389/1"
$ws.Cells.Item(139,3).Value2 = 321
$ws.Cells.Item(140,1).Value2 = 138
$ws.Cells.Item(140,2).Value = "389/10"
$ws.Cells.Item(140,3).Value2 = 321
$ws.Cells.Item(141,1).Value2 = 139
$ws.Cells.Item(141,2).Value = "390/1"
$ws.Cells.Item(141,3).Value2 = 321
$ws.Cells.Item(142,1).Value2 = 140
$ws.Cells.Item(142,2).Value = "390/9"
$ws.Cells.Item(142,3).Value2 = 321
$ws.Cells.Item(143,1).Value2 = 141
$ws.Cells.Item(143,2).Value = "390/10"
$ws.Cells.Item(143,3).Value2 = 321
$ws.Cells.Item(144,1).Value2 = 142
$ws.Cells.Item(144,2).Value = "391/1"
$ws.Cells.Item(144,3).Value2 = 321
$ws.Cells.Item(145,1).Value2 = 143
$ws.Cells.Item(145,2).Value = "391/38"
$ws.Cells.Item(145,3).Value2 = 321
$ws.Cells.Item(146,1).Value2 = 144
$ws.Cells.Item(146,2).Value = "391/45"
$ws.Cells.Item(146,3).Value2 = 321
$ws.Cells.Item(147,1).Value2 = 145
$ws.Cells.Item(147,2).Value = "392/1"
$ws.Cells.Item(147,3).Value2 = 321
$ws.Cells.Item(148,1).Value2 = 146
$ws.Cells.Item(148,2).Value = "392/6"
$ws.Cells.Item(148,3).Value2 = 321
$ws.Cells.Item(149,1).Value2 = 147
$ws.Cells.Item(149,2).Value = "392/7"
$ws.Cells.Item(149,3).Value2 = 321
$ws.Cells.Item(150,1).Value2 = 148
$ws.Cells.Item(150,2).Value = "475/21"
$ws.Cells.Item(150,3).Value2 = 321
$ws.Cells.Item(151,1).Value2 = 149
$ws.Cells.Item(151,2).Value = "475/31"
$ws.Cells.Item(151,3).Value2 = 321
$ws.Cells.Item(152,1).Value2 = 150
$ws.Cells.Item(152,2).Value = "729"
$ws.Cells.Item(152,3).Value2 = 321
$ws.Cells.Item(153,1).Value2 = 151
$ws.Cells.Item(153,2).Value = ".578"
$ws.Cells.Item(153,3).Value2 = 321
$ws.Cells.Item(154,1).Value2 = 152
$ws.Cells.Item(154,2).Value = ".579"
$ws.Cells.Item(154,3).Value2 = 321
$ws.Cells.Item(155,1).Value2 = 153
$ws.Cells.Item(155,2).Value = ".598"
$ws.Cells.Item(155,3).Value2 = 321
$ws.Cells.Item(156,1).Value2 = 154
$ws.Cells.Item(156,2).Value = ".607"
$ws.Cells.Item(156,3).Value2 = 321
$ws.Cells.Item(157,1).Value2 = 155
$ws.Cells.Item(157,2).Value = ".608"
$ws.Cells.Item(157,3).Value2 = 321
$ws.Cells.Item(158,1).Value2 = 156
$ws.Cells.Item(158,2).Value = ".451"
$ws.Cells.Item(158,3).Value2 = 88

# Column B should carry no explicit style (matches the rest of the sheet)
$bRange.ClearFormats()

# Give the new rows column-A cells the same look as the existing data rows (border + bold + centered)
$ws.Range("A114").Copy()
$ws.Range("A115:A158").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "done"